# MHV-17222: update MIN / MAX values (version + date metadata) and
# restore the "applyAlignment" formatting flag on the two body styles.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump Version and Date values -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = "0.2.10-beta"
$ws1.Range("B8").Value = "2023-12-06T12:46:33-06:00"

# --- Re-apply WrapText so the "applyAlignment" flag is written out --------
# Sheet 1 ("Metadata"): header row uses the bold style, the rest of the
# populated rows use the plain bordered style.
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A2:B16").WrapText = $true

# Sheet 2 ("Include from Ages in decades"): same two styles, but column B
# only has data starting at row 3, so avoid touching the empty B2 cell.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").WrapText = $true
$ws2.Range("A2").WrapText = $true
$ws2.Range("A3:B4").WrapText = $true
